# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 09:35"

# Armenia (row 67): refreshed case numbers
$ws.Range("B67").Value = 3538
$ws.Range("C67").Value = 146
$ws.Range("D67").Value = 1430
$ws.Range("E67").Value = 2061
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 47

# Letonia (row 101): refreshed case numbers
$ws.Range("B101").Value = 950
$ws.Range("C101").Value = 4
$ws.Range("D101").Value = 627
$ws.Range("E101").Value = 305

# Rows 192/193: Belice and Nueva Caledonia swap places (names + values)
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

# Rows 198/199: Dominica and Curazao swap places (names + values)
$ws.Range("A198").Value = "Curazao"
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "Dominica"
$ws.Range("D199").Value = 15
$ws.Range("H199").Value = 0

# Rows 215/216: San Bartolome and Sahara Occidental swap places (names only, values identical)
$ws.Range("A215").Value = "Sahara Occidental"
$ws.Range("A216").Value = "San Bartolome"
